$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 41667680
$ws.Range("I40").Value = 66667450
$ws.Range("K40").Value = 66667450
$ws.Range("M40").Value = -66667275
$ws.Range("H62").Value = 2551.7368
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2636.4375
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2636.4375
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -3884.4375
$ws.Range("H65").Value = 2551.7368
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2636.4375
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 13182.1875
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -19422.1875
$ws.Range("H125").Value = 1040.125
$ws.Range("J125").Value = 979.2
$ws.Range("L125").Value = 8812.800000000001
$ws.Range("N125").Value = -13732.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4475
$ws.Range("I88").Value = 4966.6665
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 4966.6665
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -4560.6665
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 4475
$ws.Range("I91").Value = 4966.6665
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 4966.6665
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -3562.6665
$ws.Range("N91").Value = -5808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2563.6897
$ws.Range("I86").Value = 2492.55
$ws.Range("J86").Value = 2721.7778
$ws.Range("K86").Value = 2492.55
$ws.Range("L86").Value = 2721.7778
$ws.Range("M86").Value = -1369.55
$ws.Range("N86").Value = -4967.7778
$ws.Range("H89").Value = 2563.6897
$ws.Range("I89").Value = 2492.55
$ws.Range("J89").Value = 2721.7778
$ws.Range("K89").Value = 12462.75
$ws.Range("L89").Value = 13608.889
$ws.Range("M89").Value = -6846.75
$ws.Range("N89").Value = -24840.889
$ws.Range("H94").Value = 381.41666
$ws.Range("I94").Value = 320.63635
$ws.Range("K94").Value = 320.63635
$ws.Range("M94").Value = 130.36365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H31").Value = 1370.381
$ws.Range("I31").Value = 1288.9
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1288.9
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -993.9000000000001
$ws.Range("N31").Value = -3590
$ws.Range("H34").Value = 1370.381
$ws.Range("I34").Value = 1288.9
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1288.9
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1086.9
$ws.Range("N34").Value = -3404
$ws.Range("H62").Value = 2967
$ws.Range("J62").Value = 2967
$ws.Range("L62").Value = 2967
$ws.Range("N62").Value = -4215
$ws.Range("H65").Value = 2967
$ws.Range("J65").Value = 2967
$ws.Range("L65").Value = 14835
$ws.Range("N65").Value = -21075

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 747.5455
$ws.Range("I4").Value = 393.5
$ws.Range("J4").Value = 1172.4
$ws.Range("K4").Value = 1180.5
$ws.Range("L4").Value = 3517.2
$ws.Range("M4").Value = -1068.5
$ws.Range("N4").Value = -3741.2
$ws.Range("H94").Value = 2800.1052
$ws.Range("J94").Value = 2979.647
$ws.Range("L94").Value = 8938.940999999999
$ws.Range("N94").Value = -10290.941
$ws.Range("H96").Value = 70707150
$ws.Range("J96").Value = 70707150
$ws.Range("L96").Value = 212121450
$ws.Range("N96").Value = -212125568
$ws.Range("H110").Value = 3365.5557
$ws.Range("J110").Value = 3740
$ws.Range("L110").Value = 11220
$ws.Range("N110").Value = -19400
$ws.Range("H137").Value = 27278.021
$ws.Range("I137").Value = 2279.9092
$ws.Range("J137").Value = 35365.65
$ws.Range("K137").Value = 6839.7276
$ws.Range("L137").Value = 106096.95
$ws.Range("M137").Value = -1739.7276
$ws.Range("N137").Value = -116296.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1260.75
$ws.Range("I102").Value = 1152.7727
$ws.Range("J102").Value = 1656.6666
$ws.Range("K102").Value = 1152.7727
$ws.Range("L102").Value = 1656.6666
$ws.Range("M102").Value = 469.2273
$ws.Range("N102").Value = -4900.6666
$ws.Range("H126").Value = 1953
$ws.Range("I126").Value = 1681
$ws.Range("J126").Value = 2497
$ws.Range("K126").Value = 5043
$ws.Range("L126").Value = 7491
$ws.Range("M126").Value = -2573
$ws.Range("N126").Value = -12431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 30000002
$ws.Range("J2").Value = 30000002
$ws.Range("L2").Value = 30000002
$ws.Range("N2").Value = -30000226
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 2100
$ws.Range("I40").Value = 1955.5555
$ws.Range("K40").Value = 1955.5555
$ws.Range("M40").Value = -1819.5555
$ws.Range("H122").Value = 9223.619
$ws.Range("I122").Value = 16126.5
$ws.Range("J122").Value = 4975.6924
$ws.Range("K122").Value = 48379.5
$ws.Range("L122").Value = 14927.0772
$ws.Range("M122").Value = -45929.5
$ws.Range("N122").Value = -19827.0772
$ws.Range("H132").Value = 5316.2812
$ws.Range("I132").Value = 8220.533
$ws.Range("J132").Value = 2753.7058
$ws.Range("K132").Value = 24661.599
$ws.Range("L132").Value = 8261.117400000001
$ws.Range("M132").Value = -22131.599
$ws.Range("N132").Value = -13321.1174
$ws.Range("H140").Value = 44599.832
$ws.Range("J140").Value = 44599.832
$ws.Range("L140").Value = 44599.832
$ws.Range("N140").Value = -54959.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -5224
$ws.Range("H122").Value = 1797.5454
$ws.Range("J122").Value = 1841.6666
$ws.Range("L122").Value = 5524.9998
$ws.Range("N122").Value = -10424.9998
$ws.Range("H136").Value = 19609.455
$ws.Range("I136").Value = 21237.1
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 63711.3
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -61161.3
$ws.Range("N136").Value = -15099
